# The commit swaps the East-Asian theme font ("DejaVu Sans" -> "Tahoma")
# used by the document's paragraph styles, and makes sure the complex-script
# font ("DejaVu Sans") is explicitly recorded (rather than left to inherit)
# on a few styles whose <w:rPr> previously had no <w:rFonts> at all.
#
# Word's Style.Font object exposes the four rFonts slots as:
#   Font.NameAscii    -> w:rFonts/@w:ascii
#   Font.Name         -> w:rFonts/@w:ascii and @w:hAnsi
#   Font.NameFarEast  -> w:rFonts/@w:eastAsia
#   Font.NameOther    -> w:rFonts/@w:hAnsi
#   Font.NameBi       -> w:rFonts/@w:cs

$d = $word.ActiveDocument

# docDefaults/rPrDefault, Normal and Heading all had an explicit
# w:eastAsia="DejaVu Sans" -> change it to "Tahoma".
$d.Styles("Normal").Font.NameFarEast = "Tahoma"
$d.Styles("Heading").Font.NameFarEast = "Tahoma"

# List, Caption and Index gain an explicit w:rFonts w:cs="DejaVu Sans"
# (previously they had no <w:rFonts> element at all and just inherited it).
$d.Styles("List").Font.NameBi = "DejaVu Sans"
$d.Styles("Caption").Font.NameBi = "DejaVu Sans"
$d.Styles("Index").Font.NameBi = "DejaVu Sans"
